# Updated cryptos list on Sat Jun 15 18:54:16 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row.
# Some Price values look numeric (e.g. "144.90", "0.999"); those are written
# with a temporary text NumberFormat so Excel keeps them as literal strings
# (preserving trailing zeros / exact formatting) instead of auto-converting
# them to numbers, then the cell style is restored to Normal afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.154.14'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '3.560.56'
$ws.Range('E3').Value = '  +5.21%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '607.01'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.10%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '144.90'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.81%  '
$ws.Range('D7').Value = '3.560.22'
$ws.Range('E7').Value = '  +5.30%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +4.20%  '
$ws.Range('E10').Value = '  +2.75%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '7.99'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').Value = '4.166.46'
$ws.Range('E13').Value = '  +5.28%  '
$ws.Range('E14').Value = '  +4.45%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '30.13'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +2.20%  '
$ws.Range('D16').Value = '3.559.47'
$ws.Range('E16').Value = '  +5.19%  '
$ws.Range('D17').Value = '66.255.49'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('E19').Value = '  +9.72%  '
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('E21').Value = '  +2.20%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '431.01'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +4.33%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.611'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +5.74%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '78.95'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('D25').Value = '3.703.14'
$ws.Range('E25').Value = '  +5.23%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E28').Value = '  +4.85%  '
$ws.Range('E29').Value = '  +3.71%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.13'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('E33').Value = '  +5.04%  '
$ws.Range('D34').Value = '3.555.38'
$ws.Range('E34').Value = '  +5.11%  '
$ws.Range('E35').Value = '  -2.47%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  +4.69%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '7.91'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +5.69%  '
$ws.Range('E39').Value = '  +1.97%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.02%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '169.98'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.65%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0854'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('E43').Value = '  +4.21%  '
$ws.Range('E44').Value = '  +3.32%  '
$ws.Range('E45').Value = '  +1.48%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '46.19'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('E47').Value = '  +3.38%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '25.96'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.22%  '
$ws.Range('E49').Value = '  +5.81%  '
$ws.Range('E50').Value = '  +1.74%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '23.45'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +16.60%  '
